# Atualização de bases das ligas, do dia: 04-04-2024 às 23:22
# Appends 8 new Sweden Superettan fixtures (rows 173-180) to sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Template cells (row 172 is the last existing data row) whose number
# formats / styles we want the new rows' id (col A) and date (col E)
# cells to inherit.
$styleSrcA = $ws.Range("A172")
$styleSrcE = $ws.Range("E172")

$newRows = @(
    @{ Row=173; Id=171; MatchId=7640914; Date=45388.33333333334; Home="Ostersunds FK";     Away="Trelleborgs FF";   K=1.8;   L=3.6;   M=4.2;  N=1.8;   O=3.6;   P=4.2;  Q=-0.5;  R=1.825; S=2.025; T=2.5;  U=1.875; V=1.975 },
    @{ Row=174; Id=172; MatchId=7640915; Date=45388.41666666666; Home="Sandvikens IF";     Away="Orebro SK";        K=2.625; L=3.4;   M=2.375;N=2.625; O=3.4;   P=2.375;Q=0;     R=2.025; S=1.825; T=2.75; U=1.95;  V=1.9   },
    @{ Row=175; Id=173; MatchId=7640916; Date=45388.41666666666; Home="Degerfors";         Away="GIF Sundsvall";    K=1.615; L=4;     M=5;    N=1.65;  O=4;     P=4.75; Q=-0.75; R=1.825; S=2.025; T=2.75; U=1.925; V=1.925 },
    @{ Row=176; Id=174; MatchId=7640917; Date=45388.41666666666; Home="Varbergs BoIS FC";  Away="Osters IF";        K=3.1;   L=3.4;   M=2.1;  N=3.1;   O=3.4;   P=2.1;  Q=0.25;  R=1.95;  S=1.9;   T=2.5;  U=1.825; V=2.025 },
    @{ Row=177; Id=175; MatchId=7640913; Date=45389.33333333334; Home="Utsiktens BK";      Away="Gefle IF";         K=1.55;  L=4.333; M=5.25; N=1.55;  O=4.333; P=5.25; Q=-1;    R=1.975; S=1.875; T=2.75; U=1.95;  V=1.9   },
    @{ Row=178; Id=176; MatchId=7640912; Date=45389.41666666666; Home="Skvde AIK";         Away="Landskrona BoIS";  K=3;     L=3.75;  M=2.15; N=2.8;   O=3.75;  P=2.25; Q=0.25;  R=1.85;  S=2;     T=2.75; U=1.975; V=1.875 },
    @{ Row=179; Id=177; MatchId=7642160; Date=45390.58333333334;  Home="IK Oddevold";       Away="Orgryte IS";       K=2.15;  L=3.4;   M=3.1;  N=2.15;  O=3.4;   P=3.1;  Q=-0.25; R=1.925; S=1.925; T=2.5;  U=1.825; V=2.025 },
    @{ Row=180; Id=178; MatchId=7640911; Date=45391.58333333334; Home="Helsingborg";       Away="IK Brage";         K=2.375; L=3.4;   M=2.9;  N=2.375; O=3.4;   P=2.9;  Q=-0.25; R=2.1;   S=1.775; T=2.5;  U=2;     V=1.85  }
)

foreach ($r in $newRows) {
    $row = $r.Row

    # Carry over the id column's (A) bold/border style and the date
    # column's (E) custom YYYY-MM-DD HH:MM:SS number format, same as
    # every other data row in the sheet.
    $styleSrcA.Copy($ws.Range("A$row"))
    $styleSrcE.Copy($ws.Range("E$row"))

    $ws.Cells.Item($row, 1).Value = $r.Id
    $ws.Cells.Item($row, 2).Value = $r.MatchId
    $ws.Cells.Item($row, 3).Value = "Sweden Superettan"
    $ws.Cells.Item($row, 4).Value = "Sweden Superettan"
    $ws.Cells.Item($row, 5).Value = $r.Date
    $ws.Cells.Item($row, 6).Value = $r.Home
    $ws.Cells.Item($row, 7).Value = $r.Away

    # H (FTHG) / I (FTAG) / J (FTR) intentionally left blank - these
    # fixtures have not been played yet.

    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
    $ws.Cells.Item($row, 21).Value = $r.U
    $ws.Cells.Item($row, 22).Value = $r.V

    # W (PLH) .. AA (PL_Aha) are all zero - no settled bets yet.
    $ws.Cells.Item($row, 23).Value = 0
    $ws.Cells.Item($row, 24).Value = 0
    $ws.Cells.Item($row, 25).Value = 0
    $ws.Cells.Item($row, 26).Value = 0
    $ws.Cells.Item($row, 27).Value = 0

    # AB (PL_AhOver) / AC (PL_AhUnder) intentionally left blank, same
    # as H/I/J above.
}
